# Weekly update: insert a new daily record as row 432, pushing the
# existing rows (432-478) down to (433-479).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 432 — this shifts every
# row from 432..478 down to 433..479 (matches the dimension change
# A1:R478 -> A1:R479 from the diff).
$ws.Rows(432).Insert()

# Populate the newly inserted row 432 with the new data point.
$ws.Range("A432").Value = 11
$ws.Range("B432").Value = "Vega Monumental Concepción"
$ws.Range("C432").Value = "Bíobío"
$ws.Range("D432").Value = 45212
$ws.Range("E432").Value = 8
$ws.Range("F432").Value = 100114013
$ws.Range("G432").Value = "Zanahoria"
$ws.Range("H432").Value = "Sin especificar"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 150
$ws.Range("K432").Value = 5000
$ws.Range("L432").Value = 5000
$ws.Range("M432").Value = 5000
$ws.Range("N432").Value = "$/saco 20 kilos"
$ws.Range("O432").Value = "Región de La Araucanía"
$ws.Range("P432").Value = 250
$ws.Range("Q432").Value = 20
$ws.Range("R432").Value = "Hortaliza"
